$d = $word.ActiveDocument

# wdReplace constants: 0=None,1=One,2=All ; wdFindWrap: 1=wdFindContinue
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll) | Out-Null
}

# Section A - Udaje o zamestnanci
Replace-Text "Jméno a příjmení: Martin Vader" "Jméno a příjmení: Martin Vader: Martin Vader"
Replace-Text "Osobní číslo: 2015" "Osobní číslo: 2015: 2015"
Replace-Text "Útvar / oddělení: magiologie" "Útvar / oddělení: magiologie: útvar pro magiologii"
Replace-Text "Telefon / e‑mail:" "Telefon / e‑mail: mkk@magik.cz"

# Section B - Udaje o ceste
Replace-Text "Datum a čas odjezdu: 30.11. (čas bude doplněn)" "Datum a čas odjezdu: 30.11. (čas bude doplněn): 30.11.2025"
Replace-Text "Datum a čas návratu: 5.12. (čas bude doplněn)" "Datum a čas návratu: 5.12. (čas bude doplněn): 5.12.2025"
Replace-Text "Místo konání cesty (město, adresa): Fakultní nemocnice u sv. Anny v Brně" "Místo konání cesty (město, adresa): Fakultní nemocnice u sv. Anny v Brně: Brno"
Replace-Text "Účel cesty (stručný popis):" "Účel cesty (stručný popis): školení studentů o magii"

# Section C - Zpusob dopravy
Replace-Text "C — Způsob dopravy (vyberte / označte)" "C — Způsob dopravy (vyberte / označte): hromadná doprava – vlak"

# Section D - Odhadovane naklady
Replace-Text "D — Odhadované náklady (nepovinné)" "D — Odhadované náklady (nepovinné): Odhad nákladů: 80 Kč – jízdné vlakem."

# Section E - Prohlaseni zadatele
Replace-Text "Datum:" "Datum: 27.11.2025"

# Both "Elektronický podpis:" occurrences (section E and F) get a trailing space appended.
# A single Replace-All call already updates every match in the document.
Replace-Text "Elektronický podpis:" "Elektronický podpis: "
